$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("406:408").Insert()

$ws.Cells.Item(406, 1).Value = 9
$ws.Cells.Item(406, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(406, 3).Value = "Metropolitana"
$ws.Cells.Item(406, 4).Value = 44461
$ws.Cells.Item(406, 5).Value = 13
$ws.Cells.Item(406, 6).Value = "Fruta"
$ws.Cells.Item(406, 7).Value = 100102
$ws.Cells.Item(406, 8).Value = "Cítricos"
$ws.Cells.Item(406, 9).Value = 100102005
$ws.Cells.Item(406, 10).Value = "Naranja"
$ws.Cells.Item(406, 11).Value = "Lane Late"
$ws.Cells.Item(406, 12).Value = "Primera"
$ws.Cells.Item(406, 13).Value = 260
$ws.Cells.Item(406, 14).Value = 6000
$ws.Cells.Item(406, 15).Value = 6000
$ws.Cells.Item(406, 16).Value = 6000
$ws.Cells.Item(406, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(406, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(406, 19).Value = 333
$ws.Cells.Item(406, 20).Value = 18

$ws.Cells.Item(407, 1).Value = 9
$ws.Cells.Item(407, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(407, 3).Value = "Metropolitana"
$ws.Cells.Item(407, 4).Value = 44461
$ws.Cells.Item(407, 5).Value = 13
$ws.Cells.Item(407, 6).Value = "Fruta"
$ws.Cells.Item(407, 7).Value = 100102
$ws.Cells.Item(407, 8).Value = "Cítricos"
$ws.Cells.Item(407, 9).Value = 100102005
$ws.Cells.Item(407, 10).Value = "Naranja"
$ws.Cells.Item(407, 11).Value = "Lane Late"
$ws.Cells.Item(407, 12).Value = "Segunda"
$ws.Cells.Item(407, 13).Value = 220
$ws.Cells.Item(407, 14).Value = 5000
$ws.Cells.Item(407, 15).Value = 5000
$ws.Cells.Item(407, 16).Value = 5000
$ws.Cells.Item(407, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(407, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(407, 19).Value = 278
$ws.Cells.Item(407, 20).Value = 18

$ws.Cells.Item(408, 1).Value = 9
$ws.Cells.Item(408, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(408, 3).Value = "Metropolitana"
$ws.Cells.Item(408, 4).Value = 44461
$ws.Cells.Item(408, 5).Value = 13
$ws.Cells.Item(408, 6).Value = "Fruta"
$ws.Cells.Item(408, 7).Value = 100102
$ws.Cells.Item(408, 8).Value = "Cítricos"
$ws.Cells.Item(408, 9).Value = 100102005
$ws.Cells.Item(408, 10).Value = "Naranja"
$ws.Cells.Item(408, 11).Value = "Lane Late"
$ws.Cells.Item(408, 12).Value = "Tercera"
$ws.Cells.Item(408, 13).Value = 180
$ws.Cells.Item(408, 14).Value = 4000
$ws.Cells.Item(408, 15).Value = 4000
$ws.Cells.Item(408, 16).Value = 4000
$ws.Cells.Item(408, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(408, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(408, 19).Value = 222
$ws.Cells.Item(408, 20).Value = 18
